# New crime data collected
#
# Updates the weekly CompStat report (69th Precinct) to the new reporting
# week (1/20/2025 - 1/26/2025, Volume 32 Number 4) and refreshes every
# crime-count / percent-change figure in the data table (rows 15-30) with
# the newly collected numbers.
#
# Note on text vs. number cells: some figures in this table are stored as
# literal text ("0" or "***.*" -- Excel's placeholder for an undefined /
# infinite percent change) instead of as numbers, and several cells flip
# between the two representations week over week. Plain `.Value = "0"`
# would be auto-coerced back to the number 0 by Excel, so those cells are
# written with a leading apostrophe to force text, then only the *format*
# (not the value) is copied in from a neighboring cell that already has
# the right style, via PasteSpecial(xlPasteFormats).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-TextCell($ref, $text, $formatSource) {
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($formatSource).Copy()
    $ws.Range($ref).PasteSpecial($xlPasteFormats)
}

function Set-NumericStyledCell($ref, $value, $formatSource) {
    $ws.Range($ref).Value = $value
    $ws.Range($formatSource).Copy()
    $ws.Range($ref).PasteSpecial($xlPasteFormats)
}

# --- Report header (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/20/2025  Through  1/26/2025"

# --- Row 15: Murder ---
Set-TextCell "C15" "0" "D15"
Set-TextCell "G15" "0" "D15"
Set-TextCell "H15" "***.*" "E15"

# --- Row 16: Rape ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 7
$ws.Range("L16").Value = -41.666666666666
$ws.Range("M16").Value = -46.153846153846
$ws.Range("N16").Value = -85.106382978723

# --- Row 17: Robbery ---
$ws.Range("C17").Value = 1
Set-TextCell "D17" "0" "D15"
Set-TextCell "E17" "***.*" "E15"
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 10
$ws.Range("K17").Value = 42.857142857142
$ws.Range("L17").Value = -28.571428571428
$ws.Range("M17").Value = -28.571428571428
$ws.Range("N17").Value = -23.076923076923

# --- Row 18: Fel. Assault ---
Set-TextCell "C18" "0" "D15"
Set-TextCell "D18" "0" "D15"
Set-TextCell "E18" "***.*" "E15"
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = -66.666666666666
Set-NumericStyledCell "L18" -50 "M18"
$ws.Range("M18").Value = -95.833333333333
$ws.Range("N18").Value = -97.826086956521

# --- Row 19: Burglary ---
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = -40
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 24
$ws.Range("K19").Value = -41.666666666666
$ws.Range("L19").Value = -30
$ws.Range("M19").Value = 7.692307692307
$ws.Range("N19").Value = -44

# --- Row 20: Gr. Larceny ---
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -66.666666666666
$ws.Range("J20").Value = 14
$ws.Range("K20").Value = -71.428571428571
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = -71.428571428571
$ws.Range("N20").Value = -97.530864197530

# --- Row 21: TOTAL (bold) ---
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = -18.181818181818
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 60
$ws.Range("H21").Value = -28.333333333333
$ws.Range("I21").Value = 38
$ws.Range("J21").Value = 56
$ws.Range("K21").Value = -32.142857142857
$ws.Range("L21").Value = -30.909090909090
$ws.Range("M21").Value = -51.282051282051
$ws.Range("N21").Value = -87.118644067796

# --- Row 22: Transit ---
Set-NumericStyledCell "D22" 1 "G14"
Set-NumericStyledCell "E22" -100 "H14"
$ws.Range("G22").Value = 2
$ws.Range("J22").Value = 2

# --- Row 23: Housing ---
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = -42.857142857142
$ws.Range("I23").Value = 8
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = -38.461538461538
$ws.Range("L23").Value = 33.333333333333
$ws.Range("M23").Value = 700

# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 30
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = -34.782608695652
$ws.Range("I24").Value = 23
$ws.Range("J24").Value = 42
$ws.Range("K24").Value = -45.238095238095
$ws.Range("L24").Value = -48.888888888888
$ws.Range("M24").Value = -36.111111111111

# --- Row 25: Retail Theft ---
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = -7.692307692307
$ws.Range("I25").Value = 10
$ws.Range("J25").Value = 10
$ws.Range("L25").Value = -16.666666666666

# --- Row 26: Misd. Assault ---
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -24.242424242424
$ws.Range("I26").Value = 22
$ws.Range("J26").Value = 28
$ws.Range("K26").Value = -21.428571428571
$ws.Range("L26").Value = 37.5
$ws.Range("M26").Value = -31.25

# --- Row 27: UCR Rape* ---
Set-TextCell "C27" "0" "D15"
Set-TextCell "G27" "0" "D15"
Set-TextCell "H27" "***.*" "E15"

# --- Row 28: Other Sex Crimes ---
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 3
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = -40
$ws.Range("L28").Value = 200

# --- Row 29: Shooting Vic. ---
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 200
$ws.Range("L29").Value = -50
$ws.Range("N29").Value = -66.666666666666

# --- Row 30: Shooting Inc. ---
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("L30").Value = -50
$ws.Range("N30").Value = -66.666666666666

Write-Host "Applied weekly crime data update (rows 8-9, 15-30)."
